# Add a new worksheet "add_new_departments" as the last sheet, matching
# the test-data pattern used by the other "add_new_*" sheets in this
# workbook, and make it the active/selected sheet.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "add_new_departments"

# --- Header row -----------------------------------------------------
$ws.Range("A1").Value = "code"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "head"
$ws.Range("D1").Value = "faculty"
$ws.Range("E1").Value = "status"
$ws.Range("F1").Value = "runmode"

# --- Column B (name) --------------------------------------------------
$ws.Range("B2").Value = "Exmaination"
$ws.Range("B3").Value = "IT"
$ws.Range("B4").Value = "Mathematics"
$ws.Range("B5").Value = "Graphics"
$ws.Range("B6").Value = "Networking"

# --- Column C (head) --------------------------------------------------
$ws.Range("C2").Value = "Sarath"
$ws.Range("C3").Value = "Mahin"
$ws.Range("C4").Value = "Sarath"
$ws.Range("C5").Value = "Mahin"
$ws.Range("C6").Value = "Mahin"

# --- Column A (code) ----------------------------------------------------
$ws.Range("A2").Value = "DEPT008"
$ws.Range("A3").Value = "DEPT009"
$ws.Range("A4").Value = "DEPT010"
$ws.Range("A5").Value = "DEPT011"
$ws.Range("A6").Value = "DEPT012"

# --- Columns D (faculty), E (status), F (runmode) ----------------------
$ws.Range("D2").Value = "IT"
$ws.Range("E2").Value = "Active"
$ws.Range("F2").Value = "Y"

$ws.Range("D3").Value = "IT"
$ws.Range("E3").Value = "Active"
$ws.Range("F3").Value = "Y"

$ws.Range("D4").Value = "IT"
$ws.Range("E4").Value = "Active"
$ws.Range("F4").Value = "Y"

$ws.Range("D5").Value = "IT"
$ws.Range("E5").Value = "Inactive"
$ws.Range("F5").Value = "Y"

$ws.Range("D6").Value = "IT"
$ws.Range("E6").Value = "Inactive"
$ws.Range("F6").Value = "Y"

# --- Selection / activation -------------------------------------------
$ws.Activate()
$ws.Range("A2:A6").Select()
